$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Generated PowerShell lines
$ws.Range("D2").Value = '''60.862.37'
$ws.Range("E2").Value = '  -2.43%  '
$ws.Range("D3").Value = '''2.395.93'
$ws.Range("E3").Value = '  -2.33%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  +0.14%  '
$ws.Range("D5").Value = '''569.69'
$ws.Range("E5").Value = '  -1.15%  '
$ws.Range("D6").Value = '''139.98'
$ws.Range("E6").Value = '  -2.76%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("D9").Value = '''2.395.71'
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("E10").Value = '  -0.89%  '
$ws.Range("E11").Value = '  -0.34%  '
$ws.Range("E12").Value = '  -2.27%  '
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("D14").Value = '''26.01'
$ws.Range("E14").Value = '  -1.21%  '
$ws.Range("D15").Value = '''0.0000171'
$ws.Range("E15").Value = '  -2.79%  '
$ws.Range("D16").Value = '''2.827.01'
$ws.Range("E16").Value = '  -2.39%  '
$ws.Range("D17").Value = '''60.754.85'
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("D18").Value = '''2.397.79'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").Value = '''10.65'
$ws.Range("E19").Value = '  -2.20%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").Value = '''7.42'
$ws.Range("E20").Value = '  +3.98%  '
$ws.Range("D21").Value = '''322.93'
$ws.Range("E21").Value = '  -1.77%  '
$ws.Range("E22").Value = '  -1.42%  '
$ws.Range("D23").Value = '''6.06'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = '''1.87'
$ws.Range("E25").Value = '  -4.35%  '
$ws.Range("D26").Value = '''64.91'
$ws.Range("E26").Value = '  -1.21%  '
$ws.Range("B27").Value = 'Bittensor'
$ws.Range("C27").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D27").Value = '''578.13'
$ws.Range("E27").Value = '  -2.42%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D28").Value = '''8.43'
$ws.Range("E28").Value = '  -8.63%  '
$ws.Range("D29").Value = '''2.524.54'
$ws.Range("E29").Value = '  -1.97%  '
$ws.Range("D30").Value = '''0.0₃0922'
$ws.Range("E30").Value = '  -3.90%  '
$ws.Range("E31").Value = '  -1.56%  '
$ws.Range("D32").Value = '''1.35'
$ws.Range("E32").Value = '  -6.18%  '
$ws.Range("E33").Value = '  -3.11%  '
$ws.Range("D34").Value = '''0.133'
$ws.Range("E34").Value = '  -2.78%  '
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").Value = '''4.65'
$ws.Range("E36").Value = '  -5.92%  '
$ws.Range("E37").Value = '  -3.13%  '
$ws.Range("D38").Value = '''0.367'
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("D39").Value = '''149.46'
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("E40").Value = '  -0.64%  '
$ws.Range("D41").Value = '''5.15'
$ws.Range("E41").Value = '  -3.71%  '
$ws.Range("D42").Value = '''0.999'
$ws.Range("E42").Value = '  -0.03%  '
$ws.Range("E43").Value = '  -3.75%  '
$ws.Range("D44").Value = '''41.10'
$ws.Range("E44").Value = '  -3.61%  '
$ws.Range("E45").Value = '  -4.97%  '
$ws.Range("D46").Value = '''0.0₆0284'
$ws.Range("E46").Value = '  +14.27%  '
$ws.Range("D47").Value = '''140.95'
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("E48").Value = '  -3.30%  '
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").Value = '''19.52'
$ws.Range("E50").Value = '  -1.45%  '
$ws.Range("E51").Value = '  -3.25%  '
